$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program_choosing")

# Update the program names (column A) with the new list, keep "Yes" in column B
$ws.Range("A2").Value = "TUM_MMT"
$ws.Range("A3").Value = "TUM_Consumer_Science"
$ws.Range("A4").Value = "Uni_Koeln_BA"
$ws.Range("A5").Value = "Uni_Mannheim_MGM"
$ws.Range("A6").Value = "Uni_Magdeburg_Finalcial_Economics"

# Widen column A to fit the longer program names (32 characters wide)
$ws.Columns("A").ColumnWidth = 31.17

# Move the active selection to A7
$ws.Range("A7").Select()
